$d = $word.ActiveDocument

function Get-ParagraphContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# 1) "לסיים אנימציה של שלב פתיחה" is done -> flip its highlight from yellow to
#    green (covers the run text as well as the paragraph-mark formatting).
$p1 = Get-ParagraphContaining $d "לסיים אנימציה של שלב פתיחה"
if ($p1 -ne $null) {
    $p1.Range.Font.HighlightColorIndex = 4
}

# 2) Tidy up the final "sound effect" bullet: merge the two runs
#    "להוסיף אפקט קולי" + " לאדים שיוצאים מה " into one run with the combined
#    text (the trailing "PIIPE" run is left as its own run), then mark the
#    whole line (all runs + paragraph mark) green now that it's complete too.
$rng = $d.Content
$found = $rng.Find.Execute(
    "להוסיף אפקט קולי לאדים שיוצאים מה ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "להוסיף אפקט קולי לאדים שיוצאים מה ", 2)

$p2 = Get-ParagraphContaining $d "PIIPE"
if ($p2 -ne $null) {
    $p2.Range.Font.HighlightColorIndex = 4
}
